$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5366756858871611
$ws.Range("B3").Value = 0.4690044873290731
$ws.Range("B4").Value = 0.4273327161748739
$ws.Range("B5").Value = 0.41032172068509
$ws.Range("B6").Value = 0.4074953142979325
$ws.Range("B7").Value = 0.4271034175466184
$ws.Range("B8").Value = 0.5133684956384457
$ws.Range("B9").Value = 0.6815300338288921
$ws.Range("B10").Value = 0.8044229604793145
$ws.Range("B11").Value = 0.8601793739451864
$ws.Range("B12").Value = 0.8812705911192324
$ws.Range("B13").Value = 0.8767292413427867
$ws.Range("B14").Value = 0.8619150194993495
$ws.Range("B15").Value = 0.8528379111937738
$ws.Range("B16").Value = 0.8007760616400219
$ws.Range("B17").Value = 0.7687989806938731
$ws.Range("B18").Value = 0.7503927331953264
$ws.Range("B19").Value = 0.7441583478884581
$ws.Range("B20").Value = 0.7722044405117003
$ws.Range("B21").Value = 0.8662669366005389
$ws.Range("B22").Value = 0.9276103508232154
$ws.Range("B23").Value = 0.8948826830007874
$ws.Range("B24").Value = 0.7706649011794298
$ws.Range("B25").Value = 0.6361498533019301
$ws.Range("C2").Value = 0.1457201269627433
$ws.Range("C3").Value = 0.1338437435756532
$ws.Range("C4").Value = 0.126510425852004
$ws.Range("C5").Value = 0.1235119140817318
$ws.Range("C6").Value = 0.1230134086476795
$ws.Range("C7").Value = 0.1264700275687005
$ws.Range("C8").Value = 0.1416338276000886
$ws.Range("C9").Value = 0.1710346890110941
$ws.Range("C10").Value = 0.1924214734178236
$ws.Range("C11").Value = 0.2021024881209144
$ws.Range("C12").Value = 0.2057613379665497
$ws.Range("C13").Value = 0.2049736606713282
$ws.Range("C14").Value = 0.2024036484106944
$ws.Range("C15").Value = 0.200828504512117
$ws.Range("C16").Value = 0.1917878105572868
$ws.Range("C17").Value = 0.1862291832091785
$ws.Range("C18").Value = 0.1830275100442407
$ws.Range("C19").Value = 0.181942713559323
$ws.Range("C20").Value = 0.1868213758666855
$ws.Range("C21").Value = 0.203158719056006
$ws.Range("C22").Value = 0.2137943842357117
$ws.Range("C23").Value = 0.2081218293660925
$ws.Range("C24").Value = 0.1865536640641494
$ws.Range("C25").Value = 0.1631178710582333
$ws.Range("E2").Value = 0.613646833126019
$ws.Range("E3").Value = 0.6089153485827552
$ws.Range("E4").Value = 0.6062896402596039
$ws.Range("E5").Value = 0.6052899652907726
$ws.Range("E6").Value = 0.605128218952018
$ws.Range("E7").Value = 0.6062758734884568
$ws.Range("E8").Value = 0.6119574486484112
$ws.Range("E9").Value = 0.6253151451924026
$ws.Range("E10").Value = 0.6364804907568882
$ws.Range("E11").Value = 0.6418535734021162
$ws.Range("E12").Value = 0.6439304640231427
$ws.Range("E13").Value = 0.6434812913684667
$ws.Range("E14").Value = 0.6420235945212553
$ws.Range("E15").Value = 0.6411362102351319
$ws.Range("E16").Value = 0.6361352577080339
$ws.Range("E17").Value = 0.6331425821052221
$ws.Range("E18").Value = 0.6314489410994
$ws.Range("E19").Value = 0.6308802564867406
$ws.Range("E20").Value = 0.6334582944307243
$ws.Range("E21").Value = 0.6424506099015375
$ws.Range("E22").Value = 0.6485736834414482
$ws.Range("E23").Value = 0.6452831821061409
$ws.Range("E24").Value = 0.6333154771302176
$ws.Range("E25").Value = 0.6214642384611651
$ws.Range("F2").Value = 1.821140056834466
$ws.Range("F3").Value = 1.808094619239839
$ws.Range("F4").Value = 1.801158539905998
$ws.Range("F5").Value = 1.798602004845876
$ws.Range("F6").Value = 1.798193798178715
$ws.Range("F7").Value = 1.801122968549137
$ws.Range("F8").Value = 1.8164190621643
$ws.Range("F9").Value = 1.854942296967963
$ws.Range("F10").Value = 1.888459718711189
$ws.Range("F11").Value = 1.904843686326416
$ws.Range("F12").Value = 1.911211509693047
$ws.Range("F13").Value = 1.909832808402555
$ws.Range("F14").Value = 1.905364292161167
$ws.Range("F15").Value = 1.902648501906427
$ws.Range("F16").Value = 1.887411871191389
$ws.Range("F17").Value = 1.878355913295891
$ws.Range("F18").Value = 1.873254153137481
$ws.Range("F19").Value = 1.871545156592333
$ws.Range("F20").Value = 1.879308861881697
$ws.Range("F21").Value = 1.906672363566983
$ws.Range("F22").Value = 1.925509458065221
$ws.Range("F23").Value = 1.915368467787602
$ws.Range("F24").Value = 1.878877707915592
$ws.Range("F25").Value = 1.843606171466178
$ws.Range("G2").Value = 0.2596747568516378
$ws.Range("G3").Value = 0.2580350798519291
$ws.Range("G4").Value = 0.2572759575861809
$ws.Range("G5").Value = 0.2570287149077686
$ws.Range("G6").Value = 0.2569914065924195
$ws.Range("G7").Value = 0.2572723719427614
$ws.Range("G8").Value = 0.2590578957220018
$ws.Range("G9").Value = 0.2645324182651194
$ws.Range("G10").Value = 0.2697696096455076
$ws.Range("G11").Value = 0.2724186438925642
$ws.Range("G12").Value = 0.2734602964747381
$ws.Range("G13").Value = 0.2732342418356808
$ws.Range("G14").Value = 0.2725035682214809
$ws.Range("G15").Value = 0.2720610314799075
$ws.Range("G16").Value = 0.2696018695391444
$ws.Range("G17").Value = 0.2681616688411594
$ws.Range("G18").Value = 0.2673583872114165
$ws.Range("G19").Value = 0.2670907124206678
$ws.Range("G20").Value = 0.2683123832658367
$ws.Range("G21").Value = 0.2727171378807185
$ws.Range("G22").Value = 0.2758205122283783
$ws.Range("G23").Value = 0.2741435708499012
$ws.Range("G24").Value = 0.2682441683218002
$ws.Range("G25").Value = 0.2628389145119598
$ws.Range("H2").Value = 0.4282963081503013
$ws.Range("H3").Value = 0.4310242898740029
$ws.Range("H4").Value = 0.4329313498128684
$ws.Range("H5").Value = 0.4337668449456729
$ws.Range("H6").Value = 0.4339091025634971
$ws.Range("H7").Value = 0.432942381328516
$ws.Range("H8").Value = 0.4291887484785946
$ws.Range("H9").Value = 0.4236695991612009
$ws.Range("H10").Value = 0.4207382725738569
$ws.Range("H11").Value = 0.4196489222855888
$ws.Range("H12").Value = 0.4192715318820177
$ws.Range("H13").Value = 0.4193512472702849
$ws.Range("H14").Value = 0.4196171701314029
$ws.Range("H15").Value = 0.4197846301371584
$ws.Range("H16").Value = 0.4208143781369529
$ws.Range("H17").Value = 0.4215086342704666
$ws.Range("H18").Value = 0.4219309291453897
$ws.Range("H19").Value = 0.4220778567253944
$ws.Range("H20").Value = 0.421432351343995
$ws.Range("H21").Value = 0.4195381087718317
$ws.Range("H22").Value = 0.4185048432944569
$ws.Range("H23").Value = 0.4190375776201591
$ws.Range("H24").Value = 0.4214667667208687
$ws.Range("H25").Value = 0.4249654211711373
$ws.Range("J2").Value = 0.02296014546912239
$ws.Range("J3").Value = 0.02307780131247839
$ws.Range("J4").Value = 0.02315750984052656
$ws.Range("J5").Value = 0.02319187375777343
$ws.Range("J6").Value = 0.02319769365762347
$ws.Range("J7").Value = 0.02315796565833583
$ws.Range("J8").Value = 0.02299916593230122
$ws.Range("J9").Value = 0.0227468239512767
$ws.Range("J10").Value = 0.02259718924253207
$ws.Range("J11").Value = 0.02253683019725727
$ws.Range("J12").Value = 0.02251507834587407
$ws.Range("J13").Value = 0.02251971392252905
$ws.Range("J14").Value = 0.02253501853970441
$ws.Range("J15").Value = 0.02254453682154001
$ws.Range("J16").Value = 0.02260128866489808
$ws.Range("J17").Value = 0.0226380761152889
$ws.Range("J18").Value = 0.02265996127120573
$ws.Range("J19").Value = 0.02266749601862372
$ws.Range("J20").Value = 0.02263408492006391
$ws.Range("J21").Value = 0.02253049325013023
$ws.Range("J22").Value = 0.02246922795410811
$ws.Range("J23").Value = 0.02250133863481629
$ws.Range("J24").Value = 0.02263588704779806
$ws.Range("J25").Value = 0.0228087914490871
$ws.Range("M2").Value = 0.4989552367505041
$ws.Range("M3").Value = 0.4684625176467634
$ws.Range("M4").Value = 0.4498863300257412
$ws.Range("M5").Value = 0.4423536309213958
$ws.Range("M6").Value = 0.4411050940033761
$ws.Range("M7").Value = 0.4497845900402808
$ws.Range("M8").Value = 0.4884111891386667
$ws.Range("M9").Value = 0.5653056504769296
$ws.Range("M10").Value = 0.622486502212837
$ws.Range("M11").Value = 0.6486461740458509
$ws.Range("M12").Value = 0.658573075222165
$ws.Range("M13").Value = 0.6564342215345533
$ws.Range("M14").Value = 0.6494624518957721
$ws.Range("M15").Value = 0.6451947371300548
$ws.Range("M16").Value = 0.6207798417061667
$ws.Range("M17").Value = 0.6058396298163871
$ws.Range("M18").Value = 0.5972603553480553
$ws.Range("M19").Value = 0.5943579685311562
$ws.Range("M20").Value = 0.6074286014058714
$ws.Range("M21").Value = 0.6515096687195978
$ws.Range("M22").Value = 0.6804402289569538
$ws.Range("M23").Value = 0.6649885199347807
$ws.Range("M24").Value = 0.6067101960041441
$ws.Range("M25").Value = 0.5443821503341155
$ws.Range("N2").Value = 1.00672664905661
$ws.Range("N3").Value = 1.011324605028932
$ws.Range("N4").Value = 1.014552577230695
$ws.Range("N5").Value = 1.015969954476482
$ws.Range("N6").Value = 1.016211471343979
$ws.Range("N7").Value = 1.014571279460299
$ws.Range("N8").Value = 1.008228100048271
$ws.Range("N9").Value = 0.9989947070341287
$ws.Range("N10").Value = 0.9941574104763902
$ws.Range("N11").Value = 0.9923779608909911
$ws.Range("N12").Value = 0.9917645529603618
$ws.Range("N13").Value = 0.9918939754651461
$ws.Range("N14").Value = 0.9923262849788728
$ws.Range("N15").Value = 0.992598953540778
$ws.Range("N16").Value = 0.9942821644089008
$ws.Range("N17").Value = 0.9954225188898533
$ws.Range("N18").Value = 0.9961180595812209
$ws.Range("N19").Value = 0.9963603693766743
$ws.Range("N20").Value = 0.9952970245860513
$ws.Range("N21").Value = 0.9921976660263994
$ws.Range("N22").Value = 0.9905242429578749
$ws.Range("N23").Value = 0.9913851927631896
$ws.Range("N24").Value = 0.99535363615351
$ws.Range("N25").Value = 1.001150222621824
$ws.Range("O2").Value = 1.295565746960193
$ws.Range("O3").Value = 1.297555297880223
$ws.Range("O4").Value = 1.299779117795921
$ws.Range("O5").Value = 1.300936995596842
$ws.Range("O6").Value = 1.301144448850621
$ws.Range("O7").Value = 1.299793714893269
$ws.Range("O8").Value = 1.29604350866471
$ws.Range("O9").Value = 1.296659896865407
$ws.Range("O10").Value = 1.301999782396308
$ws.Range("O11").Value = 1.30549636457863
$ws.Range("O12").Value = 1.306974371503401
$ws.Range("O13").Value = 1.306649203493606
$ws.Range("O14").Value = 1.305614873713807
$ws.Range("O15").Value = 1.305001375438565
$ws.Range("O16").Value = 1.301792789031822
$ws.Range("O17").Value = 1.300098129675433
$ws.Range("O18").Value = 1.299223853523927
$ws.Range("O19").Value = 1.298945077369723
$ws.Range("O20").Value = 1.300268130032151
$ws.Range("O21").Value = 1.305914500488655
$ws.Range("O22").Value = 1.310502132836461
$ws.Range("O23").Value = 1.307971373241713
$ws.Range("O24").Value = 1.30019096139003
$ws.Range("O25").Value = 1.29563694713886
